# Insert a new weekly record at row 13 (Jengibre, Mercado Mayorista Lo
# Valledor de Santiago). This shifts the former rows 13-46 down to 14-47
# and grows the sheet from A1:R46 to A1:R47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Insert()

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44459
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100114007
$ws.Range("G13").Value = "Jengibre"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 90
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 14000
$ws.Range("M13").Value = 12956
$ws.Range("N13").Value = "$/caja 13 kilos"
$ws.Range("O13").Value = "Perú"
$ws.Range("P13").Value = 997
$ws.Range("Q13").Value = 13
$ws.Range("R13").Value = "Hortaliza"
